$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hotfix: correct the mis-identified "units per package" value (12 -> 8)
# for the affected product rows (E3:E192).
$ws.Range("E3:E192").Value = 8

# Restore the user's final cursor position/selection after the edit.
$ws.Range("E2").Select() | Out-Null
